$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H32").Value = 16428.428
$ws.Range("I32").Value = 9999.666999999999
$ws.Range("J32").Value = 21250
$ws.Range("K32").Value = 9999.666999999999
$ws.Range("L32").Value = 21250
$ws.Range("M32").Value = -9673.666999999999
$ws.Range("N32").Value = -21902
$ws.Range("H33").Value = 31579466
$ws.Range("I33").Value = 50000390
$ws.Range("J33").Value = 736
$ws.Range("K33").Value = 50000390
$ws.Range("L33").Value = 736
$ws.Range("M33").Value = -50000161
$ws.Range("N33").Value = -1194
$ws.Range("H55").Value = 216.55556
$ws.Range("I55").Value = 281.22223
$ws.Range("J55").Value = 151.88889
$ws.Range("K55").Value = 281.22223
$ws.Range("L55").Value = 151.88889
$ws.Range("M55").Value = -67.22223000000002
$ws.Range("N55").Value = -579.8888899999999
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H113").Value = 63983.367
$ws.Range("I113").Value = 82949.21000000001
$ws.Range("K113").Value = 82949.21000000001
$ws.Range("M113").Value = -79695.21000000001
$ws.Range("H138").Value = 2908.978
$ws.Range("I138").Value = 1345.1666
$ws.Range("J138").Value = 3469.1492
$ws.Range("K138").Value = 4035.4998
$ws.Range("L138").Value = 10407.4476
$ws.Range("M138").Value = 1104.5002
$ws.Range("N138").Value = -20687.4476

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28088.428
$ws.Range("I32").Value = 16816.637
$ws.Range("K32").Value = 16816.637
$ws.Range("M32").Value = -16529.637
$ws.Range("H74").Value = 100023410
$ws.Range("I74").Value = 5674.6665
$ws.Range("K74").Value = 5674.6665
$ws.Range("M74").Value = -4800.6665
$ws.Range("H77").Value = 100023410
$ws.Range("I77").Value = 5674.6665
$ws.Range("K77").Value = 28373.3325
$ws.Range("M77").Value = -24005.3325

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 64975
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 64975
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 64975
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -66347
$ws.Range("H66").Value = 64975
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 64975
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 194925
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -201789
$ws.Range("H107").Value = 1769.0588
$ws.Range("I107").Value = 2164.5557
$ws.Range("K107").Value = 2164.5557
$ws.Range("M107").Value = -244.5556999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6727852.5
$ws.Range("I31").Value = 3166709.8
$ws.Range("J31").Value = 15630708
$ws.Range("K31").Value = 3166709.8
$ws.Range("L31").Value = 15630708
$ws.Range("M31").Value = -3166414.8
$ws.Range("N31").Value = -15631298
$ws.Range("H34").Value = 6727852.5
$ws.Range("I34").Value = 3166709.8
$ws.Range("J34").Value = 15630708
$ws.Range("K34").Value = 3166709.8
$ws.Range("L34").Value = 15630708
$ws.Range("M34").Value = -3166507.8
$ws.Range("N34").Value = -15631112
$ws.Range("H58").Value = 2363.2942
$ws.Range("I58").Value = 1926.08
$ws.Range("J58").Value = 3577.7778
$ws.Range("K58").Value = 1926.08
$ws.Range("L58").Value = 3577.7778
$ws.Range("M58").Value = -1723.08
$ws.Range("N58").Value = -3983.7778
$ws.Range("H107").Value = 959.9231
$ws.Range("I107").Value = 1124.5
$ws.Range("J107").Value = 589.625
$ws.Range("K107").Value = 1124.5
$ws.Range("L107").Value = 589.625
$ws.Range("M107").Value = 795.5
$ws.Range("N107").Value = -4429.625
$ws.Range("H132").Value = 2101.05
$ws.Range("I132").Value = 2101.05
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6303.150000000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3773.150000000001
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 2363.2942
$ws.Range("I136").Value = 1926.08
$ws.Range("J136").Value = 3577.7778
$ws.Range("K136").Value = 5778.24
$ws.Range("L136").Value = 10733.3334
$ws.Range("M136").Value = -3228.24
$ws.Range("N136").Value = -15833.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 200223.2
$ws.Range("I14").Value = 200223.2
$ws.Range("K14").Value = 600669.6000000001
$ws.Range("M14").Value = -600496.6000000001
$ws.Range("H68").Value = 4546942
$ws.Range("I68").Value = 900
$ws.Range("K68").Value = 2700
$ws.Range("M68").Value = -1889
$ws.Range("H71").Value = 4546942
$ws.Range("I71").Value = 900
$ws.Range("K71").Value = 8100
$ws.Range("M71").Value = -4044

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1321.3684
$ws.Range("I97").Value = 1266.6945
$ws.Range("K97").Value = 1266.6945
$ws.Range("M97").Value = -770.6945000000001
$ws.Range("H123").Value = 49906.535
$ws.Range("J123").Value = 53592.152
$ws.Range("L123").Value = 53592.152
$ws.Range("N123").Value = -58492.152

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2833.3333
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2833.3333
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 2833.3333
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -3423.3333
$ws.Range("H27").Value = 2833.3333
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 2833.3333
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 2833.3333
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -3047.3333
$ws.Range("H46").Value = 1559.6923
$ws.Range("J46").Value = 4426.6
$ws.Range("L46").Value = 4426.6
$ws.Range("N46").Value = -4802.6
$ws.Range("H55").Value = 1135.9231
$ws.Range("I55").Value = 1390.9
$ws.Range("K55").Value = 1390.9
$ws.Range("M55").Value = -1217.9
$ws.Range("H82").Value = 4997.3335
$ws.Range("I82").Value = 3151
$ws.Range("J82").Value = 10074.75
$ws.Range("K82").Value = 3151
$ws.Range("L82").Value = 10074.75
$ws.Range("M82").Value = -2790
$ws.Range("N82").Value = -10796.75
$ws.Range("H85").Value = 4997.3335
$ws.Range("I85").Value = 3151
$ws.Range("J85").Value = 10074.75
$ws.Range("K85").Value = 3151
$ws.Range("L85").Value = 10074.75
$ws.Range("M85").Value = -1903
$ws.Range("N85").Value = -12570.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 21000
$ws.Range("J49").Value = 21000
$ws.Range("L49").Value = 21000
$ws.Range("N49").Value = -21460
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 1327.4
$ws.Range("I107").Value = 1079.0834
$ws.Range("J107").Value = 1699.875
$ws.Range("K107").Value = 3237.2502
$ws.Range("L107").Value = 5099.625
$ws.Range("M107").Value = -1317.2502
$ws.Range("N107").Value = -8939.625
$ws.Range("H112").Value = 97000
$ws.Range("J112").Value = 97000
$ws.Range("L112").Value = 97000
$ws.Range("N112").Value = -99954
$ws.Range("H125").Value = 333381000
$ws.Range("J125").Value = 333381000
$ws.Range("L125").Value = 333381000
$ws.Range("N125").Value = -333390840
$ws.Range("H136").Value = 2558.1428
$ws.Range("I136").Value = 2303.2222
$ws.Range("J136").Value = 3017
$ws.Range("K136").Value = 6909.6666
$ws.Range("L136").Value = 9051
$ws.Range("M136").Value = -4359.6666
$ws.Range("N136").Value = -14151

